# CompStat weekly report refresh: advance the report one week
# (Volume/Number caption + "Week Covering" date range), and load the
# newly collected precinct crime-complaint figures for the week.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Header: bump the issue number and roll the reporting week forward
# ---------------------------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/8/2024  Through  1/14/2024"

# ---------------------------------------------------------------
# Helper donor cells already carrying the exact number formats we
# need to reapply after a cell's underlying type changes:
#   C14  -> General/text style used for a literal "0" placeholder
#   E14  -> General/text style used for the "***.*" placeholder
#   A21  -> bold General/text style (TOTAL row) for "0" placeholder
#   F15  -> numeric "#,##0" style (count columns)
#   M16  -> numeric "#,##0.0;-#,##0.0" style (% change columns)
# ---------------------------------------------------------------

function Set-TextCell($addr, $donor, $text) {
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

function Set-NumericCell($addr, $donor, $number) {
    $ws.Range($donor).Copy() | Out-Null
    $ws.Range($addr).PasteSpecial(-4122) | Out-Null
    $ws.Range($addr).Value = $number
}

# ---------------------------------------------------------------
# Row 15 - Rape: week count now 0
# ---------------------------------------------------------------
Set-TextCell "C15" "C14" "0"

# ---------------------------------------------------------------
# Row 16 - Robbery: counts/percentages shift
# ---------------------------------------------------------------
Set-NumericCell "D16" "F15" 1
Set-NumericCell "E16" "M16" -100
Set-TextCell    "F16" "C14" "0"
Set-NumericCell "G16" "F15" 1
Set-NumericCell "H16" "M16" -100
Set-NumericCell "J16" "F15" 1
Set-NumericCell "K16" "M16" -100

# ---------------------------------------------------------------
# Row 19 - Gr. Larceny: counts/percentages shift
# ---------------------------------------------------------------
Set-TextCell    "D19" "C14" "0"
Set-TextCell    "E19" "E14" "***.*"
$ws.Range("G19").Value = 1
Set-NumericCell "L19" "M16" -100
Set-NumericCell "M19" "M16" -100

# ---------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------
Set-TextCell "C21" "A21" "0"
$ws.Range("E21").Value = -100
$ws.Range("F21").Value = 2
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 2
$ws.Range("K21").Value = -50
$ws.Range("L21").Value = -50
$ws.Range("M21").Value = -50
$ws.Range("N21").Value = -87.5

# ---------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------
$ws.Range("F24").Value = 2
$ws.Range("G24").Value = 2
$ws.Range("H24").Value = 0

# ---------------------------------------------------------------
# Row 25 - Misd. Assault
# ---------------------------------------------------------------
Set-TextCell "C25" "C14" "0"
$ws.Range("L25").Value = -50

# ---------------------------------------------------------------
# Row 26 - UCR Rape*
# ---------------------------------------------------------------
Set-TextCell "C26" "C14" "0"

# ---------------------------------------------------------------
# Row 27 - Other Sex Crimes
# ---------------------------------------------------------------
Set-TextCell "D27" "C14" "0"
Set-TextCell "E27" "E14" "***.*"
